$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from E1 (header style) onto F1, then set header text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill in time_taken values for each data row (no special style, like column E data cells)
$ws.Range("F2").Value = "2021-10-05 13:39:32.073657"
$ws.Range("F3").Value = "2021-10-05 13:39:32.073667"
$ws.Range("F4").Value = "2021-10-05 13:39:32.073670"
$ws.Range("F5").Value = "2021-10-05 13:39:32.073673"
$ws.Range("F6").Value = "2021-10-05 13:39:32.073676"
$ws.Range("F7").Value = "2021-10-05 13:39:32.073679"
$ws.Range("F8").Value = "2021-10-05 13:39:32.073681"
$ws.Range("F9").Value = "2021-10-05 13:39:32.073684"
$ws.Range("F10").Value = "2021-10-05 13:39:32.073687"
$ws.Range("F11").Value = "2021-10-05 13:39:32.073689"
$ws.Range("F12").Value = "2021-10-05 13:39:32.073692"
$ws.Range("F13").Value = "2021-10-05 13:39:32.073694"
$ws.Range("F14").Value = "2021-10-05 13:39:32.073697"
$ws.Range("F15").Value = "2021-10-05 13:39:32.073699"
$ws.Range("F16").Value = "2021-10-05 13:39:32.073702"
$ws.Range("F17").Value = "2021-10-05 13:39:32.073704"
$ws.Range("F18").Value = "2021-10-05 13:39:32.073707"
$ws.Range("F19").Value = "2021-10-05 13:39:32.073710"
$ws.Range("F20").Value = "2021-10-05 13:39:32.073713"
$ws.Range("F21").Value = "2021-10-05 13:39:32.073715"
$ws.Range("F22").Value = "2021-10-05 13:39:32.073718"
$ws.Range("F23").Value = "2021-10-05 13:39:32.073720"
$ws.Range("F24").Value = "2021-10-05 13:39:32.073723"
$ws.Range("F25").Value = "2021-10-05 13:39:32.073726"
$ws.Range("F26").Value = "2021-10-05 13:39:32.073728"
$ws.Range("F27").Value = "2021-10-05 13:39:32.073731"
$ws.Range("F28").Value = "2021-10-05 13:39:32.073734"
$ws.Range("F29").Value = "2021-10-05 13:39:32.073736"
$ws.Range("F30").Value = "2021-10-05 13:39:32.073739"
$ws.Range("F31").Value = "2021-10-05 13:39:32.073742"
$ws.Range("F32").Value = "2021-10-05 13:39:32.073744"

$excel.CutCopyMode = 0
